# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps recorded on the handback status
# workbook to reflect the time the report was (re)generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2017-02-22 07:55:11"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2017-02-22 07:54:52"
$zhcn.Range("L2").Value = "2017-02-22 07:55:52"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2017-02-22 07:55:11"
$dede.Range("L2").Value = "2017-02-22 07:56:16"
